$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object 'object[,]' 24,9

$data[0,0] = 0.03891287221836979
$data[0,1] = 0.06009487530547201
$data[0,2] = 0.168016881724305
$data[0,3] = 0.328467677382257
$data[0,4] = 0.06429032330864373
$data[0,5] = 0.02631601038559739
$data[0,6] = 0.04458374995360299
$data[0,7] = 0.06070610586828314
$data[0,8] = 0.02625403130643734

$data[1,0] = 0.02446363992367492
$data[1,1] = 0.0761761220668528
$data[1,2] = 0.180186745027481
$data[1,3] = 0.3293526162858437
$data[1,4] = 0.06016572434712225
$data[1,5] = 0.0263166356017718
$data[1,6] = 0.04648422760032089
$data[1,7] = 0.06142056634170726
$data[1,8] = 0.02878027086378428

$data[2,0] = 0.02451039661237094
$data[2,1] = 0.07572858939447284
$data[2,2] = 0.1728005870064706
$data[2,3] = 0.3286806387199075
$data[2,4] = 0.05687111609161079
$data[2,5] = 0.02631716986182796
$data[2,6] = 0.04810714177860369
$data[2,7] = 0.06062379513123379
$data[2,8] = 0.02694665914415563

$data[3,0] = 0.02456857696647773
$data[3,1] = 0.07673622196336062
$data[3,2] = 0.1703613370406578
$data[3,3] = 0.3266481359348853
$data[3,4] = 0.0551232881493333
$data[3,5] = 0.02631744752988823
$data[3,6] = 0.0489514714181125
$data[3,7] = 0.05854520659599099
$data[3,8] = 0.02875598329775192

$data[4,0] = 0.02480416069338568
$data[4,1] = 0.07858572025094274
$data[4,2] = 0.1709690044556469
$data[4,3] = 0.3251269035504705
$data[4,4] = 0.05436098684528399
$data[4,5] = 0.02631756394633487
$data[4,6] = 0.04930533273104668
$data[4,7] = 0.05701351374056506
$data[4,8] = 0.02922663095813015

$data[5,0] = 0.03073494267917295
$data[5,1] = 0.09233583519779968
$data[5,2] = 0.1867668800946493
$data[5,3] = 0.3251215614200511
$data[5,4] = 0.05436108408954962
$data[5,5] = 0.026317517551735
$data[5,6] = 0.04916429341716718
$data[5,7] = 0.05701351374056583
$data[5,8] = 0.03019840622691315

$data[6,0] = 0.03319082467437726
$data[6,1] = 0.09597408252315069
$data[6,2] = 0.1912235097038258
$data[6,3] = 0.3248235706367357
$data[6,4] = 0.05454972473951101
$data[6,5] = 0.02631747474503368
$data[6,6] = 0.04903345714202403
$data[6,7] = 0.05672830845412419
$data[6,8] = 0.03223079160565834

$data[7,0] = 0.0264965303356203
$data[7,1] = 0.08840524714791291
$data[7,2] = 0.1884129109164361
$data[7,3] = 0.3242730900397624
$data[7,4] = 0.05726984162206547
$data[7,5] = 0.02631706669240235
$data[7,6] = 0.04779306693624821
$data[7,7] = 0.0562898783611315
$data[7,8] = 0.03354349119090147

$data[8,0] = 0.03995698985630196
$data[8,1] = 0.05904220948339548
$data[8,2] = 0.1658843795828995
$data[8,3] = 0.3254672155055917
$data[8,4] = 0.06398635345049303
$data[8,5] = 0.02631606611740528
$data[8,6] = 0.04475176292375944
$data[8,7] = 0.05772817998171039
$data[8,8] = 0.02648135536980604

$data[9,0] = 0.0829405370705361
$data[9,1] = 0.04317210931167365
$data[9,2] = 0.1388368225175758
$data[9,3] = 0.3262253576438807
$data[9,4] = 0.07038497097328139
$data[9,5] = 0.02631506104408477
$data[9,6] = 0.04169515116412279
$data[9,7] = 0.05873337873935966
$data[9,8] = 0.01815227592322571

$data[10,0] = 0.0991087437951336
$data[10,1] = 0.0454703042329198
$data[10,2] = 0.1301605875141102
$data[10,3] = 0.3256506290410212
$data[10,4] = 0.07341638064130768
$data[10,5] = 0.02631455851417541
$data[10,6] = 0.04042468654576024
$data[10,7] = 0.0582919599589014
$data[10,8] = 0.01598169884716235

$data[11,0] = 0.103581953729678
$data[11,1] = 0.04621076000510743
$data[11,2] = 0.1233370739548124
$data[11,3] = 0.3250518515592679
$data[11,4] = 0.07576275563386838
$data[11,5] = 0.02631416341624188
$data[11,6] = 0.04217480723796172
$data[11,7] = 0.05779998450094727
$data[11,8] = 0.01427213131142487

$data[12,0] = 0.09828275734984229
$data[12,1] = 0.04440923458829341
$data[12,2] = 0.1215364734027647
$data[12,3] = 0.3243168437113658
$data[12,4] = 0.07617309185885532
$data[12,5] = 0.02631409485496401
$data[12,6] = 0.04614403776522878
$data[12,7] = 0.05709245099717305
$data[12,8] = 0.01397807358611772

$data[13,0] = 0.09266683984967901
$data[13,1] = 0.04342062152097199
$data[13,2] = 0.122094369750073
$data[13,3] = 0.3232770050203363
$data[13,4] = 0.07558156087020478
$data[13,5] = 0.02631419631216331
$data[13,6] = 0.04839837716495977
$data[13,7] = 0.05604249725647667
$data[13,8] = 0.01442655130188171

$data[14,0] = 0.09045391887850908
$data[14,1] = 0.04317038526940423
$data[14,2] = 0.1230895770739566
$data[14,3] = 0.3243179100560009
$data[14,4] = 0.07558721788400533
$data[14,5] = 0.02631419342961621
$data[14,6] = 0.04916309111809585
$data[14,7] = 0.05706867955092057
$data[14,8] = 0.01440951697118928

$data[15,0] = 0.08842314426991903
$data[15,1] = 0.04304356462069656
$data[15,2] = 0.1243051484185147
$data[15,3] = 0.325668219143138
$data[15,4] = 0.07562329473731634
$data[15,5] = 0.02631418486221711
$data[15,6] = 0.04972974695711557
$data[15,7] = 0.05840241224911419
$data[15,8] = 0.01436622580686315

$data[16,0] = 0.07859052814968519
$data[16,1] = 0.04360809134277315
$data[16,2] = 0.1293625999299413
$data[16,3] = 0.3256195506119802
$data[16,4] = 0.07370366765662298
$data[16,5] = 0.02631450637913846
$data[16,6] = 0.05005883312197001
$data[16,7] = 0.05827355220341381
$data[16,8] = 0.01576953475247862

$data[17,0] = 0.07423856478264189
$data[17,1] = 0.04430331935487205
$data[17,2] = 0.1331801061903943
$data[17,3] = 0.3262748466438073
$data[17,4] = 0.0725130150174213
$data[17,5] = 0.02631470343550605
$data[17,6] = 0.04882200848657078
$data[17,7] = 0.05887082051479602
$data[17,8] = 0.01662286060421072

$data[18,0] = 0.08481770164195759
$data[18,1] = 0.043110508007674
$data[18,2] = 0.1330421932454609
$data[18,3] = 0.3261559965483951
$data[18,4] = 0.07252179172372392
$data[18,5] = 0.02631470615194003
$data[18,6] = 0.04356104875981315
$data[18,7] = 0.05875358982995313
$data[18,8] = 0.01662180173617273

$data[19,0] = 0.1111172139400722
$data[19,1] = 0.05154972840376316
$data[19,2] = 0.1253062165984322
$data[19,3] = 0.3267079395653278
$data[19,4] = 0.07562679019152091
$data[19,5] = 0.02631418263465747
$data[19,6] = 0.0391731703463965
$data[19,7] = 0.05942927314575818
$data[19,8] = 0.01435096450059466

$data[20,0] = 0.1276083585226166
$data[20,1] = 0.06238896809142695
$data[20,2] = 0.1187015274139156
$data[20,3] = 0.3281213882635706
$data[20,4] = 0.07863968021302684
$data[20,5] = 0.02631366061003151
$data[20,6] = 0.03744612562406919
$data[20,7] = 0.0609548665598042
$data[20,8] = 0.01208786899712314

$data[21,0] = 0.1391828393318206
$data[21,1] = 0.07043678979731138
$data[21,2] = 0.114146645305132
$data[21,3] = 0.3296861303983449
$data[21,4] = 0.08094485126218429
$data[21,5] = 0.02631325371953004
$data[21,6] = 0.03620749194815712
$data[21,7] = 0.06260235172373035
$data[21,8] = 0.01031947770738692

$data[22,0] = 0.1429454857424432
$data[22,1] = 0.07303436587799916
$data[22,2] = 0.1129929621476027
$data[22,3] = 0.3308203919087628
$data[22,4] = 0.08180156788500809
$data[22,5] = 0.02631310057734926
$data[22,6] = 0.03574020682224014
$data[22,7] = 0.06376250975505594
$data[22,8] = 0.009649093522931468

$data[23,0] = 0.1163996170740625
$data[23,1] = 0.05436868607465137
$data[23,2] = 0.1258474725472538
$data[23,3] = 0.3318787384074893
$data[23,4] = 0.07731650519855204
$data[23,5] = 0.02631388612666142
$data[23,6] = 0.03812769427496363
$data[23,7] = 0.06461772859928673
$data[23,8] = 0.01303546384517721

$ws.Range("B2:J25").Value = $data
